$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
$ws.Range("H98").Value = 26352574
$ws.Range("I98").Value = 9524962
$ws.Range("J98").Value = 85249210
$ws.Range("K98").Value = 9524962
$ws.Range("L98").Value = 85249210
$ws.Range("M98").Value = -9523464
$ws.Range("N98").Value = -85252206

$ws.Range("H112").Value = 1230.2632
$ws.Range("I112").Value = 716.6667
$ws.Range("J112").Value = 1467.3077
$ws.Range("K112").Value = 2150.0001
$ws.Range("L112").Value = 4401.9231
$ws.Range("M112").Value = -1042.0001
$ws.Range("N112").Value = -6617.9231

$ws.Range("H122").Value = 26352574
$ws.Range("I122").Value = 9524962
$ws.Range("J122").Value = 85249210
$ws.Range("K122").Value = 28574886
$ws.Range("L122").Value = 255747630
$ws.Range("M122").Value = -28572436
$ws.Range("N122").Value = -255752530

$ws.Range("H137").Value = 10281567
$ws.Range("I137").Value = 2315540.2
$ws.Range("J137").Value = 32921854
$ws.Range("K137").Value = 6946620.600000001
$ws.Range("L137").Value = 98765562
$ws.Range("M137").Value = -6944070.600000001
$ws.Range("N137").Value = -98770662

$ws.Range("H138").Value = 2348.28
$ws.Range("I138").Value = 1875.7333
$ws.Range("J138").Value = 3057.1
$ws.Range("K138").Value = 5627.199900000001
$ws.Range("L138").Value = 9171.299999999999
$ws.Range("M138").Value = -487.1999000000005
$ws.Range("N138").Value = -19451.3

$ws = $wb.Worksheets.Item(2)  # ARM
$ws.Range("H45").Value = 371098.56
$ws.Range("I45").Value = 417290.47
$ws.Range("J45").Value = 1563.3334
$ws.Range("K45").Value = 417290.47
$ws.Range("L45").Value = 1563.3334
$ws.Range("M45").Value = -416913.47

$ws.Range("H61").Value = 2094405.5
$ws.Range("I61").Value = 1126971.2
$ws.Range("J61").Value = 5348502.5
$ws.Range("K61").Value = 1126971.2
$ws.Range("L61").Value = 5348502.5
$ws.Range("M61").Value = -1126759.2
$ws.Range("N61").Value = -5348926.5

$ws.Range("H74").Value = 91038900
$ws.Range("I74").Value = 112699340
$ws.Range("J74").Value = 66670896
$ws.Range("K74").Value = 112699340
$ws.Range("L74").Value = 66670896
$ws.Range("M74").Value = -112698466
$ws.Range("N74").Value = -66672644

$ws.Range("H77").Value = 91038900
$ws.Range("I77").Value = 112699340
$ws.Range("J77").Value = 66670896
$ws.Range("K77").Value = 563496700
$ws.Range("L77").Value = 333354480
$ws.Range("M77").Value = -563492332
$ws.Range("N77").Value = -333363216

$ws.Range("H92").Value = 27225
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 27225
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 27225
$ws.Range("N92").Value = -32217

$ws.Range("H136").Value = 2094405.5
$ws.Range("I136").Value = 1126971.2
$ws.Range("J136").Value = 5348502.5
$ws.Range("K136").Value = 3380913.6
$ws.Range("L136").Value = 16045507.5
$ws.Range("M136").Value = -3378363.6
$ws.Range("N136").Value = -16050607.5

$ws = $wb.Worksheets.Item(3)  # BSM
$ws.Range("H107").Value = 542.5454999999999
$ws.Range("I107").Value = 617.2222
$ws.Range("J107").Value = 206.5
$ws.Range("K107").Value = 617.2222
$ws.Range("L107").Value = 206.5
$ws.Range("M107").Value = 1302.7778
$ws.Range("N107").Value = -4046.5

$ws.Range("H134").Value = 20605788
$ws.Range("I134").Value = 27779010
$ws.Range("J134").Value = 4466038.5
$ws.Range("K134").Value = 83337030
$ws.Range("L134").Value = 13398115.5
$ws.Range("M134").Value = -83334495
$ws.Range("N134").Value = -13403185.5

$ws = $wb.Worksheets.Item(4)  # CRP
$ws.Range("H31").Value = 3256901.8
$ws.Range("I31").Value = 1544809.5
$ws.Range("J31").Value = 12502200
$ws.Range("K31").Value = 1544809.5
$ws.Range("L31").Value = 12502200
$ws.Range("M31").Value = -1544514.5

$ws.Range("H34").Value = 3256901.8
$ws.Range("I34").Value = 1544809.5
$ws.Range("J34").Value = 12502200
$ws.Range("K34").Value = 1544809.5
$ws.Range("L34").Value = 12502200
$ws.Range("M34").Value = -1544607.5

$ws.Range("H58").Value = 1751223.8
$ws.Range("I58").Value = 4045.4707
$ws.Range("J58").Value = 5051449.5
$ws.Range("K58").Value = 4045.4707
$ws.Range("L58").Value = 5051449.5
$ws.Range("M58").Value = -3842.4707
$ws.Range("N58").Value = -5051855.5

$ws.Range("H132").Value = 1700.0889
$ws.Range("I132").Value = 1232.4546
$ws.Range("J132").Value = 2986.0833
$ws.Range("K132").Value = 3697.3638
$ws.Range("L132").Value = 8958.249899999999
$ws.Range("M132").Value = -1167.3638

$ws.Range("H134").Value = 913126.2
$ws.Range("I134").Value = 4063.2354
$ws.Range("J134").Value = 4003940.2
$ws.Range("K134").Value = 12189.7062
$ws.Range("L134").Value = 12011820.6
$ws.Range("M134").Value = -9654.706200000001
$ws.Range("N134").Value = -12016890.6

$ws.Range("H136").Value = 1751223.8
$ws.Range("I136").Value = 4045.4707
$ws.Range("J136").Value = 5051449.5
$ws.Range("K136").Value = 12136.4121
$ws.Range("L136").Value = 15154348.5
$ws.Range("M136").Value = -9586.4121
$ws.Range("N136").Value = -15159448.5

$ws = $wb.Worksheets.Item(5)  # CUL
$ws.Range("H63").Value = 3558.5715
$ws.Range("I63").Value = 2970
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 8910
$ws.Range("L63").Value = 12000
$ws.Range("M63").Value = -8161
$ws.Range("N63").Value = -13498

$ws.Range("H66").Value = 3558.5715
$ws.Range("I66").Value = 2970
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 26730
$ws.Range("L66").Value = 36000
$ws.Range("M66").Value = -22986
$ws.Range("N66").Value = -43488

$ws.Range("H121").Value = 2086477.2
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 2086477.2
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 6259431.6
$ws.Range("N121").Value = -6262051.6
$ws.Range("M121").ClearContents()

$ws = $wb.Worksheets.Item(6)  # GSM
$ws.Range("H102").Value = 6334.9414
$ws.Range("I102").Value = 6668.375
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 6668.375
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = -5046.375
$ws.Range("N102").Value = -4244

$ws.Range("H126").Value = 15100.714
$ws.Range("I126").Value = 15100.714
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 45302.142
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -42832.142
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item(7)  # LTW
$ws.Range("H61").Value = 1775.381
$ws.Range("I61").Value = 1175.2727
$ws.Range("J61").Value = 2435.5
$ws.Range("K61").Value = 1175.2727
$ws.Range("L61").Value = 2435.5
$ws.Range("M61").Value = -973.2727
$ws.Range("N61").Value = -2839.5

$ws.Range("H82").Value = 3886.4783
$ws.Range("I82").Value = 1331.2667
$ws.Range("J82").Value = 8677.5
$ws.Range("K82").Value = 1331.2667
$ws.Range("L82").Value = 8677.5
$ws.Range("M82").Value = -970.2666999999999
$ws.Range("N82").Value = -9399.5

$ws.Range("H85").Value = 3886.4783
$ws.Range("I85").Value = 1331.2667
$ws.Range("J85").Value = 8677.5
$ws.Range("K85").Value = 1331.2667
$ws.Range("L85").Value = 8677.5
$ws.Range("M85").Value = -83.2666999999999
$ws.Range("N85").Value = -11173.5

$ws.Range("H100").Value = 5882.2
$ws.Range("I100").Value = 1575
$ws.Range("J100").Value = 7448.4546
$ws.Range("K100").Value = 1575
$ws.Range("L100").Value = 7448.4546
$ws.Range("M100").Value = -1034

$ws.Range("H113").Value = 1775.381
$ws.Range("I113").Value = 1175.2727
$ws.Range("J113").Value = 2435.5
$ws.Range("K113").Value = 1175.2727
$ws.Range("L113").Value = 2435.5
$ws.Range("M113").Value = 994.7273
$ws.Range("N113").Value = -6775.5

$ws.Range("H132").Value = 17876766
$ws.Range("I132").Value = 47653380
$ws.Range("J132").Value = 10799.6
$ws.Range("K132").Value = 142960140
$ws.Range("L132").Value = 32398.8
$ws.Range("M132").Value = -142957610
$ws.Range("N132").Value = -37458.8

$ws.Range("H136").Value = 1784142.5
$ws.Range("I136").Value = 2180040.8
$ws.Range("J136").Value = 2600.6667
$ws.Range("K136").Value = 6540122.399999999
$ws.Range("L136").Value = 7802.000100000001
$ws.Range("M136").Value = -6537572.399999999
$ws.Range("N136").Value = -12902.0001

$ws = $wb.Worksheets.Item(8)  # WVR
$ws.Range("H101").Value = 14000
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 14000
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 14000
$ws.Range("N101").Value = -20490

$ws.Range("H126").Value = 22729124
$ws.Range("I126").Value = 35714564
$ws.Range("J126").Value = 4601.25
$ws.Range("K126").Value = 107143692
$ws.Range("L126").Value = 13803.75
$ws.Range("M126").Value = -107141222
$ws.Range("N126").Value = -18743.75

$ws.Range("H132").Value = 893929
$ws.Range("I132").Value = 2625.1304
$ws.Range("J132").Value = 2757564.2
$ws.Range("K132").Value = 7875.3912
$ws.Range("L132").Value = 8272692.600000001
$ws.Range("M132").Value = -5345.3912
$ws.Range("N132").Value = -8277752.600000001

$ws.Range("H136").Value = 10940.131
$ws.Range("I136").Value = 7209.0557
$ws.Range("J136").Value = 24372
$ws.Range("K136").Value = 21627.1671
$ws.Range("L136").Value = 73116
$ws.Range("M136").Value = -19077.1671
$ws.Range("N136").Value = -78216
